$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (date moved from 06-07 to 06-08)
$ws.Name = "Through 2022-06-08"

# Update the "June (through 06-07)" label cell to "June (through 06-08)"
$ws.Range("A7").Value = "June (through 06-08)"

# Update June row (row 7) values
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 13
$ws.Range("E7").Value = 18
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 40
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 29

# Update Total row (row 8) values
$ws.Range("C8").Value = 221
$ws.Range("D8").Value = 329
$ws.Range("E8").Value = 313
$ws.Range("F8").Value = 211
$ws.Range("G8").Value = 398
$ws.Range("H8").Value = 661
$ws.Range("I8").Value = 692
